$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Top Gainers")

$ws.Cells.Item(2, 3).Value = 11.551
$ws.Cells.Item(2, 4).Value = 20.387
$ws.Cells.Item(2, 5).Value = 27.5689
$ws.Cells.Item(3, 3).Value = 11.2388
$ws.Cells.Item(3, 4).Value = 8.2645
$ws.Cells.Item(3, 5).Value = -6.703
$ws.Cells.Item(4, 2).Value = "MCLOUD"
$ws.Cells.Item(4, 3).Value = 11.0036
$ws.Cells.Item(4, 4).Value = 9.9396
$ws.Cells.Item(4, 5).Value = -19.5689
$ws.Cells.Item(5, 2).Value = "SAGILITY"
$ws.Cells.Item(5, 3).Value = 10.635
$ws.Cells.Item(5, 4).Value = 18.4343
$ws.Cells.Item(5, 5).Value = 31.9269
$ws.Cells.Item(6, 3).Value = 10.5298
$ws.Cells.Item(6, 4).Value = 10.7364
$ws.Cells.Item(6, 5).Value = 24.8056
$ws.Cells.Item(7, 3).Value = 10.0508
$ws.Cells.Item(7, 4).Value = 15.3192
$ws.Cells.Item(7, 5).Value = 22.2097
$ws.Cells.Item(8, 2).Value = "MARINE"
$ws.Cells.Item(8, 3).Value = 9.1496
$ws.Cells.Item(8, 4).Value = 5.587
$ws.Cells.Item(8, 5).Value = 18.4808
$ws.Cells.Item(9, 2).Value = "UNIPARTS"
$ws.Cells.Item(9, 3).Value = 8.427899999999999
$ws.Cells.Item(9, 4).Value = 10.72
$ws.Cells.Item(9, 5).Value = 27.022
$ws.Cells.Item(10, 2).Value = "CHENNPETRO"
$ws.Cells.Item(10, 3).Value = 7.5584
$ws.Cells.Item(10, 4).Value = 12.7598
$ws.Cells.Item(10, 5).Value = 14.9248
$ws.Cells.Item(11, 3).Value = 6.4392
$ws.Cells.Item(11, 4).Value = 6.7256
$ws.Cells.Item(11, 5).Value = 8.476599999999999
$ws.Cells.Item(12, 3).Value = 6.2042
$ws.Cells.Item(12, 4).Value = 10.4823
$ws.Cells.Item(12, 5).Value = 11.5781
$ws.Cells.Item(13, 3).Value = 5.6937
$ws.Cells.Item(13, 4).Value = 19.0875
$ws.Cells.Item(13, 5).Value = 19.1763
$ws.Cells.Item(14, 3).Value = 5.6517
$ws.Cells.Item(14, 4).Value = 8.999000000000001
$ws.Cells.Item(14, 5).Value = 15.1658
$ws.Cells.Item(15, 2).Value = "BLSE"
$ws.Cells.Item(15, 3).Value = 5.3798
$ws.Cells.Item(15, 4).Value = 4.3463
$ws.Cells.Item(15, 5).Value = -1.845
$ws.Cells.Item(16, 2).Value = "VENKEYS"
$ws.Cells.Item(16, 3).Value = 5.3278
$ws.Cells.Item(16, 4).Value = 5.9566
$ws.Cells.Item(16, 5).Value = 3.652
$ws.Cells.Item(17, 2).Value = "BAJAJHCARE"
$ws.Cells.Item(17, 3).Value = 5.3044
$ws.Cells.Item(17, 4).Value = 5.8579
$ws.Cells.Item(17, 5).Value = -0.494
$ws.Cells.Item(18, 3).Value = 5.1184
$ws.Cells.Item(18, 4).Value = 11.6957
$ws.Cells.Item(18, 5).Value = 8.1778
$ws.Cells.Item(19, 2).Value = "V2RETAIL"
$ws.Cells.Item(19, 3).Value = 4.9975
$ws.Cells.Item(19, 4).Value = 3.4793
$ws.Cells.Item(19, 5).Value = 11.7421
$ws.Cells.Item(20, 2).Value = "SKMEGGPROD"
$ws.Cells.Item(20, 3).Value = 4.9909
$ws.Cells.Item(20, 4).Value = 12.0154
$ws.Cells.Item(20, 5).Value = 29.9408
$ws.Cells.Item(21, 2).Value = "NETWEB"
$ws.Cells.Item(21, 3).Value = 4.9157
$ws.Cells.Item(21, 4).Value = 10.7611
$ws.Cells.Item(21, 5).Value = 12.8436
$ws.Cells.Item(22, 2).Value = "POLICYBZR"
$ws.Cells.Item(22, 3).Value = 4.8915
$ws.Cells.Item(22, 4).Value = 7.2373
$ws.Cells.Item(22, 5).Value = 6.2103
$ws.Cells.Item(23, 2).Value = "CREDITACC"
$ws.Cells.Item(23, 3).Value = 4.8657
$ws.Cells.Item(23, 4).Value = 3.4798
$ws.Cells.Item(23, 5).Value = 8.779199999999999
$ws.Cells.Item(24, 2).Value = "MRPL"
$ws.Cells.Item(24, 3).Value = 4.8614
$ws.Cells.Item(24, 4).Value = 15.0437
$ws.Cells.Item(24, 5).Value = 25.8905
$ws.Cells.Item(25, 3).Value = 4.738
$ws.Cells.Item(25, 4).Value = 5.2468
$ws.Cells.Item(25, 5).Value = 4.7785
$ws.Cells.Item(26, 3).Value = 4.5755
$ws.Cells.Item(26, 4).Value = 4.7545
$ws.Cells.Item(26, 5).Value = -0.3609
$ws.Cells.Item(27, 2).Value = "BGRENERGY"
$ws.Cells.Item(27, 3).Value = 4.3291
$ws.Cells.Item(27, 4).Value = -5.0783
$ws.Cells.Item(27, 5).Value = 76.71469999999999
$ws.Cells.Item(28, 2).Value = "DEEDEV"
$ws.Cells.Item(28, 3).Value = 4.3126
$ws.Cells.Item(28, 4).Value = -2.6271
$ws.Cells.Item(28, 5).Value = -3.4302
$ws.Cells.Item(29, 2).Value = "EUROPRATIK"
$ws.Cells.Item(29, 3).Value = 4.3105
$ws.Cells.Item(29, 4).Value = 10.1829
$ws.Cells.Item(29, 5).Value = 27.1375
$ws.Cells.Item(30, 2).Value = "SHANTIGOLD"
$ws.Cells.Item(30, 3).Value = 4.0436
$ws.Cells.Item(30, 4).Value = 11.3666
$ws.Cells.Item(30, 5).Value = 3.9366
$ws.Cells.Item(31, 2).Value = "ABREL"
$ws.Cells.Item(31, 3).Value = 3.8398
$ws.Cells.Item(31, 4).Value = 12.1969
$ws.Cells.Item(31, 5).Value = 11.732
$ws.Cells.Item(32, 2).Value = "HIRECT"
$ws.Cells.Item(32, 3).Value = 3.8249
$ws.Cells.Item(32, 4).Value = 11.2539
$ws.Cells.Item(32, 5).Value = 9.545500000000001
$ws.Cells.Item(33, 2).Value = "IIFL"
$ws.Cells.Item(33, 3).Value = 3.8241
$ws.Cells.Item(33, 4).Value = 10.8616
$ws.Cells.Item(33, 5).Value = 20.1593
$ws.Cells.Item(34, 2).Value = "BEML"
$ws.Cells.Item(34, 3).Value = 3.7217
$ws.Cells.Item(34, 4).Value = 0.7864
$ws.Cells.Item(34, 5).Value = 7.1685
$ws.Cells.Item(35, 2).Value = "BLS"
$ws.Cells.Item(35, 3).Value = 3.6587
$ws.Cells.Item(35, 4).Value = 0.597
$ws.Cells.Item(35, 5).Value = -0.6671
$ws.Cells.Item(36, 2).Value = "SOLEX"
$ws.Cells.Item(36, 3).Value = 3.651
$ws.Cells.Item(36, 4).Value = 5.2022
$ws.Cells.Item(36, 5).Value = "N/A"
$ws.Cells.Item(37, 2).Value = "SAPPHIRE"
$ws.Cells.Item(37, 3).Value = 3.6198
$ws.Cells.Item(37, 4).Value = 5.4469
$ws.Cells.Item(37, 5).Value = 2.791
$ws.Cells.Item(38, 2).Value = "RAMASTEEL"
$ws.Cells.Item(38, 3).Value = 3.6145
$ws.Cells.Item(38, 4).Value = 3.5105
$ws.Cells.Item(38, 5).Value = 5.0916
$ws.Cells.Item(39, 2).Value = "SKYGOLD"
$ws.Cells.Item(39, 3).Value = 3.6046
$ws.Cells.Item(39, 4).Value = -0.9437
$ws.Cells.Item(39, 5).Value = 37.5688
$ws.Cells.Item(40, 3).Value = 3.5553
$ws.Cells.Item(40, 4).Value = 3.8021
$ws.Cells.Item(40, 5).Value = 5.232
$ws.Cells.Item(41, 2).Value = "MTARTECH"
$ws.Cells.Item(41, 3).Value = 3.5317
$ws.Cells.Item(41, 4).Value = 7.6539
$ws.Cells.Item(41, 5).Value = 31.5815
$ws.Cells.Item(42, 2).Value = "PSPPROJECT"
$ws.Cells.Item(42, 3).Value = 3.4645
$ws.Cells.Item(42, 4).Value = 17.6193
$ws.Cells.Item(42, 5).Value = 24.0859
$ws.Cells.Item(44, 2).Value = "CENTRUM"
$ws.Cells.Item(44, 3).Value = 3.3333
$ws.Cells.Item(44, 4).Value = 1.9432
$ws.Cells.Item(44, 5).Value = 1.2771
$ws.Cells.Item(46, 2).Value = "RSYSTEMS"
$ws.Cells.Item(46, 3).Value = 3.2336
$ws.Cells.Item(46, 4).Value = 4.3924
$ws.Cells.Item(46, 5).Value = 6.7365
$ws.Cells.Item(47, 2).Value = "BPCL"
$ws.Cells.Item(47, 3).Value = 3.2318
$ws.Cells.Item(47, 4).Value = 8.7456
$ws.Cells.Item(47, 5).Value = 5.8001
$ws.Cells.Item(48, 2).Value = "JKTYRE"
$ws.Cells.Item(48, 3).Value = 3.1239
$ws.Cells.Item(48, 4).Value = 6.1779
$ws.Cells.Item(48, 5).Value = 22.2983
$ws.Cells.Item(49, 2).Value = "NBCC"
$ws.Cells.Item(49, 3).Value = 3.0637
$ws.Cells.Item(49, 4).Value = 6.3211
$ws.Cells.Item(49, 5).Value = 10.8984
$ws.Cells.Item(50, 2).Value = "CENTUM"
$ws.Cells.Item(50, 3).Value = 3.0057
$ws.Cells.Item(50, 4).Value = 3.8656
$ws.Cells.Item(50, 5).Value = -1.6411
$ws.Cells.Item(51, 2).Value = "SUNDROP"
$ws.Cells.Item(51, 3).Value = 2.9851
$ws.Cells.Item(51, 4).Value = 2.8778
$ws.Cells.Item(51, 5).Value = 0.9778
$ws.Cells.Item(52, 2).Value = "ALICON"
$ws.Cells.Item(52, 3).Value = 2.9765
$ws.Cells.Item(52, 4).Value = 9.081
$ws.Cells.Item(52, 5).Value = 14.4994
$ws.Cells.Item(53, 2).Value = "POWERINDIA"
$ws.Cells.Item(53, 3).Value = 2.883
$ws.Cells.Item(53, 4).Value = 7.2702
$ws.Cells.Item(53, 5).Value = -0.0833
$ws.Cells.Item(54, 2).Value = "GANESHCP"
$ws.Cells.Item(54, 3).Value = 2.8683
$ws.Cells.Item(54, 4).Value = 2.3428
$ws.Cells.Item(54, 5).Value = 1.8941
$ws.Cells.Item(55, 3).Value = 2.8366
$ws.Cells.Item(55, 4).Value = -1.5229
$ws.Cells.Item(55, 5).Value = 8.4765
$ws.Cells.Item(57, 2).Value = "DBCORP"
$ws.Cells.Item(57, 3).Value = 2.7678
$ws.Cells.Item(57, 4).Value = 5.4075
$ws.Cells.Item(57, 5).Value = 1.3556
$ws.Cells.Item(58, 2).Value = "GREAVESCOT"
$ws.Cells.Item(58, 3).Value = 2.7083
$ws.Cells.Item(58, 4).Value = 15.2825
$ws.Cells.Item(58, 5).Value = 11.5223
$ws.Cells.Item(59, 2).Value = "AHLUCONT"
$ws.Cells.Item(59, 3).Value = 2.7082
$ws.Cells.Item(59, 4).Value = 1.5196
$ws.Cells.Item(59, 5).Value = -5.6823
$ws.Cells.Item(60, 2).Value = "CARYSIL"
$ws.Cells.Item(60, 3).Value = 2.6999
$ws.Cells.Item(60, 4).Value = 2.184
$ws.Cells.Item(60, 5).Value = 11.0748
$ws.Cells.Item(61, 2).Value = "IVALUE"
$ws.Cells.Item(61, 3).Value = 2.6836
$ws.Cells.Item(61, 4).Value = 5.9208
$ws.Cells.Item(61, 5).Value = -1.2545
$ws.Cells.Item(62, 3).Value = 2.6789
$ws.Cells.Item(62, 4).Value = 3.4943
$ws.Cells.Item(62, 5).Value = 11.2003
$ws.Cells.Item(63, 2).Value = "GRAPHITE"
$ws.Cells.Item(63, 3).Value = 2.6485
$ws.Cells.Item(63, 4).Value = 16.3177
$ws.Cells.Item(63, 5).Value = 16.5271
$ws.Cells.Item(64, 2).Value = "DIVISLAB"
$ws.Cells.Item(64, 3).Value = 2.6259
$ws.Cells.Item(64, 4).Value = 1.3651
$ws.Cells.Item(64, 5).Value = 17.462
$ws.Cells.Item(65, 2).Value = "REFEX"
$ws.Cells.Item(65, 3).Value = 2.5952
$ws.Cells.Item(65, 4).Value = 0.0275
$ws.Cells.Item(65, 5).Value = 2.0053
$ws.Cells.Item(66, 2).Value = "CIFL"
$ws.Cells.Item(66, 3).Value = 2.588
$ws.Cells.Item(66, 4).Value = 2.1424
$ws.Cells.Item(66, 5).Value = 2.0833
$ws.Cells.Item(67, 2).Value = "HEG"
$ws.Cells.Item(67, 3).Value = 2.5749
$ws.Cells.Item(67, 4).Value = 15.4838
$ws.Cells.Item(67, 5).Value = 17.954
$ws.Cells.Item(68, 2).Value = "WEBELSOLAR"
$ws.Cells.Item(68, 3).Value = 2.5138
$ws.Cells.Item(68, 4).Value = 2.9105
$ws.Cells.Item(68, 5).Value = -0.9553
$ws.Cells.Item(69, 2).Value = "MAMATA"
$ws.Cells.Item(69, 3).Value = 2.5028
$ws.Cells.Item(69, 4).Value = 1.9142
$ws.Cells.Item(69, 5).Value = 1.1172
$ws.Cells.Item(70, 2).Value = "DBL"
$ws.Cells.Item(70, 3).Value = 2.4863
$ws.Cells.Item(70, 4).Value = 3.5699
$ws.Cells.Item(70, 5).Value = 4.6435
$ws.Cells.Item(72, 2).Value = "GENUSPOWER"
$ws.Cells.Item(72, 3).Value = 2.4368
$ws.Cells.Item(72, 4).Value = 10.3872
$ws.Cells.Item(72, 5).Value = 7.1118
$ws.Cells.Item(73, 2).Value = "REDTAPE"
$ws.Cells.Item(73, 3).Value = 2.4208
$ws.Cells.Item(73, 4).Value = 2.3214
$ws.Cells.Item(73, 5).Value = -4.5271
$ws.Cells.Item(74, 2).Value = "JSFB"
$ws.Cells.Item(74, 3).Value = 2.3919
$ws.Cells.Item(74, 4).Value = 2.2541
$ws.Cells.Item(74, 5).Value = -4.6833
$ws.Cells.Item(75, 2).Value = "CAMS"
$ws.Cells.Item(75, 3).Value = 2.3862
$ws.Cells.Item(75, 4).Value = 1.7344
$ws.Cells.Item(75, 5).Value = 5.0258
$ws.Cells.Item(76, 2).Value = "JKLAKSHMI"
$ws.Cells.Item(76, 3).Value = 2.3437
$ws.Cells.Item(76, 4).Value = 4.328
$ws.Cells.Item(76, 5).Value = 1.3529
